$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New data rows to append (mirrors rows that betexplorer scraper would add)
# ---------------------------------------------------------------------------
$newRows = @(
  @{
    A=63; B="kuwait"; C="premier-league"; D="2023-2024"; E=45295.64930555555;
    F="Al-Fahaheel"; G=2; H="Al Jahra"; I=0;
    J=2.17; K="04/01/2024 03:42"; L=1.93; M="04/01/2024 15:33";
    N=3.22; O="04/01/2024 03:42"; P=3.43; Q="04/01/2024 15:33";
    R=3.11; S="04/01/2024 03:42"; T=3.6;  U="04/01/2024 15:32";
    V="https://www.betexplorer.com/football/kuwait/premier-league/al-fahaheel-al-jahra/xKl5RIhb/"
  },
  @{
    A=64; B="kuwait"; C="premier-league"; D="2023-2024"; E=45295.75;
    F="Kazma SC"; G=0; H="Al Salmiya"; I=0;
    J=2.46; K="04/01/2024 06:11"; L=2.19; M="04/01/2024 17:59";
    N=3.34; O="04/01/2024 06:11"; P=3.65; Q="04/01/2024 17:59";
    R=2.58; S="04/01/2024 06:11"; T=2.82; U="04/01/2024 17:59";
    V="https://www.betexplorer.com/football/kuwait/premier-league/kazma-sc-al-salmiya/G40dTvOo/"
  }
)

$lastRow = 63
$rowIndex = $lastRow + 1

foreach ($row in $newRows) {
    # Copy the formatting (styles) of the previous data row onto the new one,
    # then overwrite the values so the new row looks/behaves like the others.
    $srcRange = $ws.Range("A" + ($rowIndex - 1) + ":V" + ($rowIndex - 1))
    $dstRange = $ws.Range("A" + $rowIndex + ":V" + $rowIndex)
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)

    $ws.Cells.Item($rowIndex, 1).Value = $row.A
    $ws.Cells.Item($rowIndex, 2).Value = $row.B
    $ws.Cells.Item($rowIndex, 3).Value = $row.C
    $ws.Cells.Item($rowIndex, 4).Value = $row.D
    $ws.Cells.Item($rowIndex, 5).Value = $row.E
    $ws.Cells.Item($rowIndex, 6).Value = $row.F
    $ws.Cells.Item($rowIndex, 7).Value = $row.G
    $ws.Cells.Item($rowIndex, 8).Value = $row.H
    $ws.Cells.Item($rowIndex, 9).Value = $row.I
    $ws.Cells.Item($rowIndex, 10).Value = $row.J
    $ws.Cells.Item($rowIndex, 11).Value = $row.K
    $ws.Cells.Item($rowIndex, 12).Value = $row.L
    $ws.Cells.Item($rowIndex, 13).Value = $row.M
    $ws.Cells.Item($rowIndex, 14).Value = $row.N
    $ws.Cells.Item($rowIndex, 15).Value = $row.O
    $ws.Cells.Item($rowIndex, 16).Value = $row.P
    $ws.Cells.Item($rowIndex, 17).Value = $row.Q
    $ws.Cells.Item($rowIndex, 18).Value = $row.R
    $ws.Cells.Item($rowIndex, 19).Value = $row.S
    $ws.Cells.Item($rowIndex, 20).Value = $row.T
    $ws.Cells.Item($rowIndex, 21).Value = $row.U
    $ws.Cells.Item($rowIndex, 22).Value = $row.V

    $rowIndex = $rowIndex + 1
}
